$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 77: correct the debit amount that was entered ---
$ws.Range("B77").Value = 46800

# D77 previously carried a stray italic style; bring it in line with the
# rest of the "Ordered Amount" column by copying the format from D76.
$ws.Range("D76").Copy($ws.Range("D77"))

# --- Row 78: new ledger entry (27-Feb-2020) ---
$ws.Range("A78").Value = 43888
$ws.Range("B78").Value = 33280
$ws.Range("D78").Value = "Ordered Amount"

# --- Row 79: new ledger entry (28-Feb-2020) ---
$ws.Range("A79").Value = 43889
$ws.Range("B79").Value = 40560
$ws.Range("D79").Value = "Ordered Amount"

# The running-balance formulas in column E already span this range as a
# shared formula; force each one to refresh its cached value now that the
# rows above it have real data (the recalc engine doesn't always ripple a
# shared-formula chain automatically).
$ws.Range("E77").Formula = $ws.Range("E77").Formula
$ws.Range("E78").Formula = $ws.Range("E78").Formula
$ws.Range("E79").Formula = $ws.Range("E79").Formula

# Leave the selection where the user left off after keying in the new
# "Ordered Amount" cells.
$ws.Range("D78:D79").Select() | Out-Null
